# Auto-generated Excel COM-interop script to apply crypto price/volume updates
# and the row-content swaps (14<->15, 29<->30, 31<->32) per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.179.41"
$ws.Range("E2").Value = "  +2.91%  "

$ws.Range("D3").Value = "2.648.74"
$ws.Range("E3").Value = "  +2.60%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'596.35"
$ws.Range("E5").Value = "  +1.02%  "

$ws.Range("D6").Value = "'156.60"
$ws.Range("E6").Value = "  +4.08%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").Value = "'0.593"
$ws.Range("E8").Value = "  +1.14%  "

$ws.Range("E9").Value = "  +7.82%  "

$ws.Range("D10").Value = "'0.400"
$ws.Range("E10").Value = "  +4.11%  "

$ws.Range("D11").Value = "'5.80"
$ws.Range("E11").Value = "  +1.67%  "

$ws.Range("E12").Value = "  +1.93%  "

$ws.Range("D13").Value = "'29.05"
$ws.Range("E13").Value = "  +5.30%  "

$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "'0.0000184"
$ws.Range("E14").Value = "  +18.36%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.121.16"
$ws.Range("E15").Value = "  +2.51%  "

$ws.Range("D16").Value = "65.075.23"
$ws.Range("E16").Value = "  +3.08%  "

$ws.Range("D17").Value = "2.731.79"
$ws.Range("E17").Value = "  +4.81%  "

$ws.Range("D18").Value = "'12.64"
$ws.Range("E18").Value = "  +3.33%  "

$ws.Range("D19").Value = "'4.82"
$ws.Range("E19").Value = "  +1.61%  "

$ws.Range("D20").Value = "'354.80"
$ws.Range("E20").Value = "  +2.64%  "

$ws.Range("D21").Value = "'7.30"
$ws.Range("E21").Value = "  +6.20%  "

$ws.Range("E22").Value = "  +0.23%  "

$ws.Range("D23").Value = "'68.23"
$ws.Range("E23").Value = "  +1.36%  "

$ws.Range("D24").Value = "'1.71"
$ws.Range("E24").Value = "  +1.23%  "

$ws.Range("D25").Value = "'9.54"
$ws.Range("E25").Value = "  +2.88%  "

$ws.Range("D26").Value = "'1.65"
$ws.Range("E26").Value = "  -1.64%  "

$ws.Range("D27").Value = "'8.17"
$ws.Range("E27").Value = "  +1.32%  "

$ws.Range("E28").Value = "  +1.08%  "

$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0949"
$ws.Range("E29").Value = "  +11.89%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'2.10"
$ws.Range("E31").Value = "  +3.43%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'519.39"
$ws.Range("E32").Value = "  -8.11%  "

$ws.Range("D33").Value = "'1.79"
$ws.Range("E33").Value = "  +1.51%  "

$ws.Range("D34").Value = "'5.65"
$ws.Range("E34").Value = "  +8.29%  "

$ws.Range("D35").Value = "'6.34"
$ws.Range("E35").Value = "  +3.80%  "

$ws.Range("D36").Value = "'0.428"
$ws.Range("E36").Value = "  +3.87%  "

$ws.Range("D37").Value = "'164.84"
$ws.Range("E37").Value = "  -1.10%  "

$ws.Range("E38").Value = "  +5.53%  "

$ws.Range("D39").Value = "'20.30"
$ws.Range("E39").Value = "  +4.24%  "

$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "  -0.14%  "

$ws.Range("E41").Value = "  +0.07%  "

$ws.Range("D42").Value = "'42.23"
$ws.Range("E42").Value = "  +6.80%  "

$ws.Range("D43").Value = "'165.53"
$ws.Range("E43").Value = "  -0.38%  "

$ws.Range("D44").Value = "'4.10"
$ws.Range("E44").Value = "  +3.24%  "

$ws.Range("D45").Value = "'0.0618"
$ws.Range("E45").Value = "  +6.21%  "

$ws.Range("D46").Value = "'23.01"
$ws.Range("E46").Value = "  +0.94%  "

$ws.Range("D47").Value = "'2.22"
$ws.Range("E47").Value = "  +7.10%  "

$ws.Range("D48").Value = "'0.650"
$ws.Range("E48").Value = "  +3.47%  "

$ws.Range("D49").Value = "'0.0256"
$ws.Range("E49").Value = "  +1.68%  "

$ws.Range("D50").Value = "'0.0986"
$ws.Range("E50").Value = "  +2.60%  "

$ws.Range("D51").Value = "'19.44"
$ws.Range("E51").Value = "  +1.55%  "

# Reset style on text-forced numeric-looking cells to drop the auto-added
# quote-prefix style and keep them on the workbook default style (no explicit s=).
$textCells = @("D4","D5","D6","D7","D8","D10","D11","D13","D14","D18","D19","D20","D21","D23","D24","D25","D26","D27","D30","D31","D32","D33","D34","D35","D36","D37","D39","D40","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
